$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Status" column (B) for the classification / turma items ---
# Row 11: "Classificar aluno e turma por ano" -> ok
$ws.Range("B11").Value = "ok"
# Row 12: "Modificar regra de lista de alunos ao criar uma turma" -> Em andamento
$ws.Range("B12").Value = "Em andamento"
# Row 14: "Emitir aviso de erro de conexao quando estiver inserindo nota" -> ok
$ws.Range("B14").Value = "ok"
# Row 20: "Alterar a id da turma do aluno tbm na tabela de boletim..." -> ok
$ws.Range("B20").Value = "ok"

# --- Give the whole status column (B10:B26) the centered / bordered look ---
# (matches the formatting already used for the "ok" cells higher up the sheet)
$rng = $ws.Range("B10:B26")
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108

# --- Move the viewport / selection like the author left it ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("G20").Select()
